$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(556).Insert()

$ws.Range("A556").Value = 5
$ws.Range("B556").Value = "Macroferia Regional de Talca"
$ws.Range("C556").Value = "Maule"
$ws.Range("D556").Value = 44984
$ws.Range("E556").Value = 7
$ws.Range("F556").Value = "Fruta"
$ws.Range("G556").Value = 100103
$ws.Range("H556").Value = "Frutos de hueso (carozo)"
$ws.Range("I556").Value = 100103004
$ws.Range("J556").Value = "Durazno"
$ws.Range("K556").Value = "White Lady"
$ws.Range("L556").Value = "Especial"
$ws.Range("M556").Value = 220
$ws.Range("N556").Value = 14000
$ws.Range("O556").Value = 14000
$ws.Range("P556").Value = 14000
$ws.Range("Q556").Value = "$/bandeja 15 kilos granel"
$ws.Range("R556").Value = "Región de O'Higgins"
$ws.Range("S556").Value = 933
$ws.Range("T556").Value = 15
